$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'52.281.61"
$ws.Range('E2').Value = "'  +1.28%  "
$ws.Range('D3').Value = "'2.903.52"
$ws.Range('E3').Value = "'  +3.94%  "
$ws.Range('E4').Value = "'  +0.14%  "
$ws.Range('D5').Value = "'352.92"
$ws.Range('E5').Value = "'  -0.12%  "
$ws.Range('D6').Value = "'113.61"
$ws.Range('E6').Value = "'  +2.04%  "
$ws.Range('D7').Value = "'0.560"
$ws.Range('E7').Value = "'  +0.86%  "
$ws.Range('E8').Value = "'  +0.08%  "
$ws.Range('D9').Value = "'0.625"
$ws.Range('E9').Value = "'  -0.57%  "
$ws.Range('D10').Value = "'40.08"
$ws.Range('E10').Value = "'  +0.03%  "
$ws.Range('D11').Value = "'0.0866"
$ws.Range('E11').Value = "'  +3.59%  "
$ws.Range('E12').Value = "'  +0.59%  "
$ws.Range('D13').Value = "'19.83"
$ws.Range('E13').Value = "'  -0.48%  "
$ws.Range('D14').Value = "'7.77"
$ws.Range('E14').Value = "'  +0.43%  "
$ws.Range('D15').Value = "'3.363.01"
$ws.Range('E15').Value = "'  +4.10%  "
$ws.Range('D16').Value = "'0.998"
$ws.Range('E16').Value = "'  +5.83%  "
$ws.Range('D17').Value = "'2.903.30"
$ws.Range('E17').Value = "'  +3.61%  "
$ws.Range('D18').Value = "'52.320.96"
$ws.Range('E18').Value = "'  +1.42%  "
$ws.Range('D19').Value = "'7.65"
$ws.Range('E19').Value = "'  +0.94%  "
$ws.Range('D20').Value = "'3.30"
$ws.Range('E20').Value = "'  +2.79%  "
$ws.Range('E21').Value = "'  +4.30%  "
$ws.Range('D22').Value = "'0.0₃0978"
$ws.Range('E22').Value = "'  +0.83%  "
$ws.Range('D23').Value = "'70.79"
$ws.Range('E23').Value = "'  +0.79%  "
$ws.Range('D24').Value = "'269.38"
$ws.Range('E24').Value = "'  +0.88%  "
$ws.Range('D25').Value = "'2.80"
$ws.Range('E25').Value = "'  +1.73%  "
$ws.Range('E26').Value = "'  +8.18%  "
$ws.Range('D27').Value = "'26.80"
$ws.Range('E27').Value = "'  +2.83%  "
$ws.Range('E28').Value = "'  -0.09%  "
$ws.Range('B29').Value = "'Hedera"
$ws.Range('C29').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('D29').Value = "'0.104"
$ws.Range('E29').Value = "'  +16.69%  "
$ws.Range('B30').Value = "'Cosmos"
$ws.Range('C30').Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range('D30').Value = "'10.61"
$ws.Range('E30').Value = "'  +2.82%  "
$ws.Range('D31').Value = "'37.64"
$ws.Range('E31').Value = "'  -3.29%  "
$ws.Range('D32').Value = "'6.57"
$ws.Range('E32').Value = "'  +7.57%  "
$ws.Range('D33').Value = "'6.25"
$ws.Range('E33').Value = "'  +12.05%  "
$ws.Range('D34').Value = "'53.20"
$ws.Range('E34').Value = "'  +1.19%  "
$ws.Range('D35').Value = "'0.0450"
$ws.Range('E35').Value = "'  -0.68%  "
$ws.Range('E36').Value = "'  -12.79%  "
$ws.Range('E37').Value = "'  +0.03%  "
$ws.Range('D38').Value = "'3.34"
$ws.Range('E38').Value = "'  +5.81%  "
$ws.Range('D39').Value = "'18.92"
$ws.Range('E39').Value = "'  +0.69%  "
$ws.Range('E40').Value = "'  +2.62%  "
$ws.Range('D41').Value = "'2.76"
$ws.Range('E41').Value = "'  +11.03%  "
$ws.Range('E42').Value = "'  +1.87%  "
$ws.Range('D43').Value = "'23.01"
$ws.Range('E43').Value = "'  +5.79%  "
$ws.Range('B44').Value = "'ApeXProtocol"
$ws.Range('C44').Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range('D44').Value = "'2.61"
$ws.Range('E44').Value = "'  +6.48%  "
$ws.Range('B45').Value = "'Monero"
$ws.Range('C45').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('D45').Value = "'120.07"
$ws.Range('E45').Value = "'  -0.70%  "
$ws.Range('E46').Value = "'  -1.97%  "
$ws.Range('E47').Value = "'  +3.83%  "
$ws.Range('D48').Value = "'2.183.18"
$ws.Range('E48').Value = "'  +3.58%  "
$ws.Range('E49').Value = "'  +21.05%  "
$ws.Range('E50').Value = "'  +13.68%  "
$ws.Range('E51').Value = "'  -0.38%  "
